$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells that would otherwise be auto-parsed as numbers,
# so the updated values are stored as text (matching the original cell type).
$textCells = @("D5", "D8", "D15", "D16", "D18", "D19", "D23", "D25", "D26", "D29", "D31", "D34", "D36", "D40", "D42", "D43", "D48", "D50")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "27.730.06"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.614.70"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "210.67"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "22.85"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "1.846.71"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "1.611.48"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = "0.550"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").Value = "64.33"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").Value = "27.752.86"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "225.58"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").Value = "7.54"
$ws.Range("E19").Value = "  -1.36%  "
$ws.Range("D20").Value = "0.0₃0709"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").Value = "9.94"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").Value = "154.81"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "6.86"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("D29").Value = "15.25"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "0.0476"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "1.393.14"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "3.04"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").Value = "0.969"
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "0.838"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").Value = "0.995"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").Value = "64.95"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("E45").Value = "  -3.78%  "
$ws.Range("D46").Value = "1.755.30"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "89.25"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").Value = "0.0987"
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("E51").Value = "  -0.55%  "

# Restore original (default) cell formatting now that the text values are set.
foreach ($c in $textCells) { $ws.Range($c).ClearFormats() }
